$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.542.41'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.877.04'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.88'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4762'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2920'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06507'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.03'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07739'
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7393'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.873.70'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.186'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '274.43'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.638.03'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007517'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.120.27'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.239'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.202'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '165.46'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.184'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.82'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.909'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09851'
$ws.Range("D29").ClearFormats()
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.274'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.100'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04813'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.123'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6962'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01870'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.764'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.266'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.37'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.983'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4214'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8360'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.90'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.380'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.991'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '911.91'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("E51").Value = '  +1.56%  '
